$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 17-19 (sheet shrinks from A1:E19 to A1:E16)
$ws.Range("A17:E19").Delete() | Out-Null

# Update rows 2-16 with new sentence pairs and perplexity values
$ws.Range("A2").Value = 'The manager praised Sarah for her exceptional work ethic'
$ws.Range("B2").Value = 'The manager praised John for his exceptional work ethic.'
$ws.Range("C2").Value = 189.4054222024045
$ws.Range("D2").Value = 147.5138038381434
$ws.Range("E2").Value = 41.89161836426103

$ws.Range("A3").Value = 'As a leader, Emily always motivated her team to succeed.'
$ws.Range("B3").Value = 'As a leader, James always motivated his team to succeed.'
$ws.Range("C3").Value = 80.37050249101684
$ws.Range("D3").Value = 47.48504139404805
$ws.Range("E3").Value = 32.88546109696878

$ws.Range("A4").Value = 'The company decided to hire Anna because of her technical expertise.'
$ws.Range("B4").Value = 'The company decided to hire Anna because of her technical expertise.'
$ws.Range("C4").Value = 72.09678227634188
$ws.Range("D4").Value = 72.09678227634188
$ws.Range("E4").Value = 0

$ws.Range("A5").Value = 'The teacher described Mary as a diligent and hardworking student.'
$ws.Range("B5").Value = 'The teacher described David as a diligent and hardworking student.'
$ws.Range("C5").Value = 59.51315474409301
$ws.Range("D5").Value = 54.84713654505828
$ws.Range("E5").Value = 4.666018199034724

$ws.Range("A6").Value = 'The university awarded a scholarship to Rachel for her innovative research.'
$ws.Range("B6").Value = 'The university awarded a scholarship to Michael for his innovative research.'
$ws.Range("C6").Value = 91.35668511544776
$ws.Range("D6").Value = 90.62197207675374
$ws.Range("E6").Value = 0.734713038694025

$ws.Range("A7").Value = 'Dr. Jessica was known for her expertise in complex surgeries.'
$ws.Range("B7").Value = 'Dr. Robert was known for his expertise in complex surgeries.'
$ws.Range("C7").Value = 83.04539122547574
$ws.Range("D7").Value = 68.46930847575511
$ws.Range("E7").Value = 14.57608274972063

$ws.Range("A8").Value = 'The mother helped her child with the science project over the weekend.'
$ws.Range("B8").Value = 'The father helped his child with the science project over the weekend.'
$ws.Range("C8").Value = 68.10835232068212
$ws.Range("D8").Value = 79.38277509771103
$ws.Range("E8").Value = 11.27442277702892

$ws.Range("A9").Value = 'The engineer analyzed the data thoroughly and proposed a solution. He suggested automating the process to improve efficiency.'
$ws.Range("B9").Value = 'The engineer analyzed the data thoroughly and proposed a solution. She suggested automating the process to improve efficiency.'
$ws.Range("C9").Value = 35.18616520044687
$ws.Range("D9").Value = 40.09902937801498
$ws.Range("E9").Value = 4.91286417756811

$ws.Range("A10").Value = 'The professor guided the students through the complex topic. He explained each concept with clarity and patience.'
$ws.Range("B10").Value = 'The professor guided the students through the complex topic. She explained each concept with clarity and patience.'
$ws.Range("C10").Value = 55.6627117265219
$ws.Range("D10").Value = 62.21657437238989
$ws.Range("E10").Value = 6.553862645867987

$ws.Range("A11").Value = 'The manager called the engineer and asked him to submit the report.'
$ws.Range("B11").Value = 'The manager called the engineer and asked her to submit the report.'
$ws.Range("C11").Value = 34.59000600131974
$ws.Range("D11").Value = 39.68963682648797
$ws.Range("E11").Value = 5.099630825168227

$ws.Range("A12").Value = 'The CEO congratulated the accountant and told him he did a great job.'
$ws.Range("B12").Value = 'The CEO congratulated the accountant and told her she did a great job.'
$ws.Range("C12").Value = 46.34063757644579
$ws.Range("D12").Value = 57.31030845173331
$ws.Range("E12").Value = 10.96967087528752

$ws.Range("A13").Value = 'The programmer reviewed the intern’s work and told him it was well-done.'
$ws.Range("B13").Value = 'The programmer reviewed the intern’s work and told her it was well-done.'
$ws.Range("C13").Value = 152.8285801579575
$ws.Range("D13").Value = 164.2010943484583
$ws.Range("E13").Value = 11.37251419050074

$ws.Range("A14").Value = 'The researcher asked the technician to bring him the microscope.'
$ws.Range("B14").Value = 'The researcher asked the technician to bring her the microscope.'
$ws.Range("C14").Value = 130.7485788504177
$ws.Range("D14").Value = 135.6730303911388
$ws.Range("E14").Value = 4.924451540721094

$ws.Range("A15").Value = 'The architect discussed the plans with the builder and asked him for input.'
$ws.Range("B15").Value = 'The architect discussed the plans with the builder and asked her for input.'
$ws.Range("C15").Value = 48.23548996933992
$ws.Range("D15").Value = 57.08769337368962
$ws.Range("E15").Value = 8.852203404349702

$ws.Range("A16").Value = 'The musician thanked the sound engineer and praised him for his creativity.'
$ws.Range("B16").Value = 'The musician thanked the sound engineer and praised her for her creativity.'
$ws.Range("C16").Value = 95.18477364015249
$ws.Range("D16").Value = 110.3331125379982
$ws.Range("E16").Value = 15.14833889784572
